# refactor to version 1.1.6
# Rearranges columns B/C/D on the analysis sheet: the numeric "count" column
# moves from D to C, the tag-list column moves from C to D, and the
# example-list column moves from B to C (i.e. B<-old C, C<-old D, D<-old B).
# Row 6 (I-NP) additionally gets refreshed example/tag contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; B = "A, Np, FW"; C = 3953; D = "nhiều, hơn, nhất, gần, khác, cùng, lớn, cao, đầu tiên, nghèo" },
    @{ Row = 3; B = "N, P, M, Np, Nc, L, Nu, FW, Ny, V"; C = 26980; D = "một, người, ông, những, các, năm, tôi, khi, nhà, anh" },
    @{ Row = 4; B = "E, Np"; C = 3817; D = "của, trong, với, ở, cho, để, đến, vào, từ, trên" },
    @{ Row = 5; B = "V, Vy, FW"; C = 12634; D = "là, có, đi, được, làm, phải, lên, ra, bị, nói" },
    @{ Row = 6; B = "Np, N, M, CH, V, Ny, A, FW, C, X"; C = 1892; D = "Văn, Nam, Hoà, Thị, Sang, HCM, Cảnh, Oanh, Việt, Á" },
    @{ Row = 7; B = "CH, R, C, T, X, I, Z, FW"; C = 16821; D = ",, ., “, ”, và, đã, không, :, ..., -" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
